# Rename the two existing sheets.
$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item(1).Name = "fund"
$wb.Worksheets.Item(2).Name = "share"

# Update the remembered selection on the "fund" sheet (was C84, now C193).
[void]$wb.Worksheets.Item("fund").Range("C193").Select()

# Add a brand new worksheet "count" at the end of the workbook.
$countSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$countSheet.Name = "count"

# Populate it with a single header cell.
$countSheet.Range("A1").Value = "时间"

# Match the original author's print setup on the new sheet.
$countSheet.PageSetup.PaperSize = 9
$countSheet.PageSetup.Orientation = 1

# Leave the selection on A2 and make "count" the active/visible tab,
# exactly like in the committed workbook.
[void]$countSheet.Range("A2").Select()
